$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 878.1818
$ws.Cells.Item(43, 9).Value = 733.3333
$ws.Cells.Item(43, 10).Value = 932.5
$ws.Cells.Item(43, 11).Value = 733.3333
$ws.Cells.Item(43, 12).Value = 932.5
$ws.Cells.Item(43, 13).Value = -664.3333
$ws.Cells.Item(43, 14).Value = -1070.5
$ws.Cells.Item(64, 8).Value = 3017.5278
$ws.Cells.Item(64, 9).Value = 2669.1428
$ws.Cells.Item(64, 10).Value = 3505.2666
$ws.Cells.Item(64, 11).Value = 2669.1428
$ws.Cells.Item(64, 12).Value = 3505.2666
$ws.Cells.Item(64, 13).Value = -2421.1428
$ws.Cells.Item(64, 14).Value = -4001.2666
$ws.Cells.Item(67, 8).Value = 3017.5278
$ws.Cells.Item(67, 9).Value = 2669.1428
$ws.Cells.Item(67, 10).Value = 3505.2666
$ws.Cells.Item(67, 11).Value = 2669.1428
$ws.Cells.Item(67, 12).Value = 3505.2666
$ws.Cells.Item(67, 13).Value = -1811.1428
$ws.Cells.Item(67, 14).Value = -5221.2666
$ws.Cells.Item(86, 8).Value = 1993.7333
$ws.Cells.Item(86, 9).Value = 1800
$ws.Cells.Item(86, 10).Value = 2122.889
$ws.Cells.Item(86, 11).Value = 1800
$ws.Cells.Item(86, 12).Value = 2122.889
$ws.Cells.Item(86, 13).Value = -677
$ws.Cells.Item(86, 14).Value = -4368.889
$ws.Cells.Item(87, 8).Value = 60000
$ws.Cells.Item(87, 10).Value = 60000
$ws.Cells.Item(87, 12).Value = 60000
$ws.Cells.Item(87, 14).Value = -62496
$ws.Cells.Item(89, 8).Value = 1993.7333
$ws.Cells.Item(89, 9).Value = 1800
$ws.Cells.Item(89, 10).Value = 2122.889
$ws.Cells.Item(89, 11).Value = 9000
$ws.Cells.Item(89, 12).Value = 10614.445
$ws.Cells.Item(89, 13).Value = -3384
$ws.Cells.Item(89, 14).Value = -21846.445
$ws.Cells.Item(90, 8).Value = 60000
$ws.Cells.Item(90, 10).Value = 60000
$ws.Cells.Item(90, 12).Value = 180000
$ws.Cells.Item(90, 14).Value = -192480
$ws.Cells.Item(112, 8).Value = 1845.1333
$ws.Cells.Item(112, 10).Value = 1845.1333
$ws.Cells.Item(112, 12).Value = 5535.3999
$ws.Cells.Item(112, 14).Value = -7751.3999
$ws.Cells.Item(115, 8).Value = 9091332
$ws.Cells.Item(115, 9).Value = 9091332
$ws.Cells.Item(115, 11).Value = 27273996
$ws.Cells.Item(115, 13).Value = -27272429
$ws.Cells.Item(135, 8).Value = 579.1316
$ws.Cells.Item(135, 9).Value = 420.80768
$ws.Cells.Item(135, 10).Value = 922.1667
$ws.Cells.Item(135, 11).Value = 3787.26912
$ws.Cells.Item(135, 12).Value = 8299.5003
$ws.Cells.Item(135, 13).Value = -1252.26912
$ws.Cells.Item(135, 14).Value = -13369.5003
$ws.Cells.Item(137, 8).Value = 18248
$ws.Cells.Item(137, 9).Value = 1015.69446
$ws.Cells.Item(137, 10).Value = 42108.117
$ws.Cells.Item(137, 11).Value = 3047.08338
$ws.Cells.Item(137, 12).Value = 126324.351
$ws.Cells.Item(137, 13).Value = -497.08338
$ws.Cells.Item(137, 14).Value = -131424.351
$ws.Cells.Item(138, 8).Value = 1958.3273
$ws.Cells.Item(138, 9).Value = 1223
$ws.Cells.Item(138, 10).Value = 2778.5
$ws.Cells.Item(138, 11).Value = 3669
$ws.Cells.Item(138, 12).Value = 8335.5
$ws.Cells.Item(138, 13).Value = 1471
$ws.Cells.Item(138, 14).Value = -18615.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 10643878
$ws.Cells.Item(32, 9).Value = 12050783
$ws.Cells.Item(32, 11).Value = 12050783
$ws.Cells.Item(32, 13).Value = -12050496
$ws.Cells.Item(92, 8).Value = 20055.715
$ws.Cells.Item(92, 10).Value = 20055.715
$ws.Cells.Item(92, 12).Value = 20055.715
$ws.Cells.Item(92, 14).Value = -25047.715
$ws.Cells.Item(132, 8).Value = 1679.6608
$ws.Cells.Item(132, 9).Value = 1267.575
$ws.Cells.Item(132, 10).Value = 2709.875
$ws.Cells.Item(132, 11).Value = 3802.725
$ws.Cells.Item(132, 12).Value = 8129.625
$ws.Cells.Item(132, 13).Value = -1272.725
$ws.Cells.Item(132, 14).Value = -13189.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2754.8
$ws.Cells.Item(105, 9).Value = 1634.4
$ws.Cells.Item(105, 10).Value = 2964.875
$ws.Cells.Item(105, 11).Value = 1634.4
$ws.Cells.Item(105, 12).Value = 2964.875
$ws.Cells.Item(105, 13).Value = 112.5999999999999
$ws.Cells.Item(105, 14).Value = -6458.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 5590.9375
$ws.Cells.Item(31, 9).Value = 0
$ws.Cells.Item(31, 10).Value = 5590.9375
$ws.Cells.Item(31, 11).Value = 0
$ws.Cells.Item(31, 12).Value = 5590.9375
$ws.Cells.Item(31, 13).ClearContents()
$ws.Cells.Item(31, 14).Value = -6180.9375
$ws.Cells.Item(34, 8).Value = 5590.9375
$ws.Cells.Item(34, 9).Value = 0
$ws.Cells.Item(34, 10).Value = 5590.9375
$ws.Cells.Item(34, 11).Value = 0
$ws.Cells.Item(34, 12).Value = 5590.9375
$ws.Cells.Item(34, 13).ClearContents()
$ws.Cells.Item(34, 14).Value = -5994.9375
$ws.Cells.Item(58, 8).Value = 1719.3096
$ws.Cells.Item(58, 9).Value = 1029.0385
$ws.Cells.Item(58, 10).Value = 2841
$ws.Cells.Item(58, 11).Value = 1029.0385
$ws.Cells.Item(58, 12).Value = 2841
$ws.Cells.Item(58, 13).Value = -826.0385000000001
$ws.Cells.Item(58, 14).Value = -3247
$ws.Cells.Item(134, 8).Value = 1683.0282
$ws.Cells.Item(134, 9).Value = 1099.8545
$ws.Cells.Item(134, 10).Value = 3687.6875
$ws.Cells.Item(134, 11).Value = 3299.5635
$ws.Cells.Item(134, 12).Value = 11063.0625
$ws.Cells.Item(134, 13).Value = -764.5634999999997
$ws.Cells.Item(134, 14).Value = -16133.0625
$ws.Cells.Item(136, 8).Value = 1719.3096
$ws.Cells.Item(136, 9).Value = 1029.0385
$ws.Cells.Item(136, 10).Value = 2841
$ws.Cells.Item(136, 11).Value = 3087.1155
$ws.Cells.Item(136, 12).Value = 8523
$ws.Cells.Item(136, 13).Value = -537.1155000000003
$ws.Cells.Item(136, 14).Value = -13623

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 4310897
$ws.Cells.Item(113, 9).Value = 8621190
$ws.Cells.Item(113, 10).Value = 603.75
$ws.Cells.Item(113, 11).Value = 25863570
$ws.Cells.Item(113, 12).Value = 1811.25
$ws.Cells.Item(113, 13).Value = -25861400
$ws.Cells.Item(113, 14).Value = -6151.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3007
$ws.Cells.Item(80, 9).Value = 2701.7368
$ws.Cells.Item(80, 10).Value = 3732
$ws.Cells.Item(80, 11).Value = 2701.7368
$ws.Cells.Item(80, 12).Value = 3732
$ws.Cells.Item(80, 13).Value = -1703.7368
$ws.Cells.Item(80, 14).Value = -5728
$ws.Cells.Item(83, 8).Value = 3007
$ws.Cells.Item(83, 9).Value = 2701.7368
$ws.Cells.Item(83, 10).Value = 3732
$ws.Cells.Item(83, 11).Value = 13508.684
$ws.Cells.Item(83, 12).Value = 18660
$ws.Cells.Item(83, 13).Value = -8516.684000000001
$ws.Cells.Item(83, 14).Value = -28644

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value = 22135
$ws.Cells.Item(98, 10).Value = 22135
$ws.Cells.Item(98, 12).Value = 22135
$ws.Cells.Item(98, 14).Value = -28125
$ws.Cells.Item(113, 8).Value = 267.9643
$ws.Cells.Item(113, 9).Value = 316.63635
$ws.Cells.Item(113, 10).Value = 89.5
$ws.Cells.Item(113, 11).Value = 949.90905
$ws.Cells.Item(113, 12).Value = 268.5
$ws.Cells.Item(113, 13).Value = 1220.09095
$ws.Cells.Item(113, 14).Value = -4608.5
$ws.Cells.Item(136, 8).Value = 2501.4783
$ws.Cells.Item(136, 9).Value = 2451.6
$ws.Cells.Item(136, 10).Value = 2697.4285
$ws.Cells.Item(136, 11).Value = 7354.799999999999
$ws.Cells.Item(136, 12).Value = 8092.2855
$ws.Cells.Item(136, 13).Value = -4804.799999999999
$ws.Cells.Item(136, 14).Value = -13192.2855
